$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial value that is the same (46075)
# for every data row (2 through 234). Bump it by one day to 46076.
$ws.Range("C2:C234").Value = 46076
